$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap dates between rows 4-5 and rows 6-7
$ws.Range("D4").Value = 44574
$ws.Range("D5").Value = 44574
$ws.Range("D6").Value = 44559
$ws.Range("D7").Value = 44559
